# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Copy the "closing" bottom-border formatting (currently on the
#    very last table row, 31) onto what will become the new last
#    data row (29) once the two trailing rows are removed.
# ------------------------------------------------------------------
$ws.Range("B31:J31").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Remove the two trailing rows belonging to GUSTAVO ADOLFO GARCIA
#    HEREDIA (periods 2103 / 2102) - this worker is no longer part
#    of the statement. Deleting shifts everything below up by 2.
# ------------------------------------------------------------------
$ws.Range("B30:J31").Delete()

# ------------------------------------------------------------------
# 3) Refresh the summary header fields.
# ------------------------------------------------------------------
$ws.Range("E11").Value2 = 637440
$ws.Range("C13").Value2 = 2
$ws.Range("F13").Value2 = 7

# ------------------------------------------------------------------
# 4) Rewrite the worker detail table (rows 16-29): two workers
#    (JESICA RODRIGUEZ TEHERAN / LESLY PATRICIA SALCEDO SAMPAYO),
#    periods 2305..2311, alternating row by row.
# ------------------------------------------------------------------
$data = @(
    @{Row=16; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2305"; Mora=46400; Salario=1000000},
    @{Row=17; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2305"; Mora=46400; Salario=1160000},
    @{Row=18; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2306"; Mora=46400; Salario=1000000},
    @{Row=19; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2306"; Mora=46400; Salario=1160000},
    @{Row=20; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2307"; Mora=46400; Salario=1000000},
    @{Row=21; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2307"; Mora=46400; Salario=1160000},
    @{Row=22; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2308"; Mora=46400; Salario=1000000},
    @{Row=23; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2308"; Mora=46400; Salario=1160000},
    @{Row=24; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2309"; Mora=46400; Salario=1000000},
    @{Row=25; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2309"; Mora=46400; Salario=1160000},
    @{Row=26; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2310"; Mora=46400; Salario=1000000},
    @{Row=27; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2310"; Mora=46400; Salario=1160000},
    @{Row=28; Doc="1128056659"; Nombre="JESICA RODRIGUEZ TEHERAN";        Periodo="2311"; Mora=37333; Salario=1000000},
    @{Row=29; Doc="45563613";   Nombre="LESLY PATRICIA SALCEDO SAMPAYO";  Periodo="2311"; Mora=43307; Salario=1160000}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value2 = $item.Mora
    $ws.Cells.Item($r, 7).Value2 = $item.Salario
}
